# feat: add 2022-Q4 data
#
# The workbook has a "总计" (totals) summary sheet plus one sheet per
# reporting quarter (newest quarter second, i.e. right after "总计").
# Adding a new quarter means:
#   1. Duplicating the current newest quarter sheet ("2022-Q3") so its
#      existing figures are preserved under their own tab.
#   2. Renaming the original tab to the new quarter ("2022-Q4") and
#      overwriting its figures with the new quarter's numbers.
#   3. Renaming the duplicate back to the quarter it actually holds
#      data for ("2022-Q3"), leaving it positioned between "2022-Q4"
#      and "2021-Q3".
#   4. Updating the "总计" sheet: the newest row becomes "2022-Q4" with
#      the new totals, and a new row is inserted beneath it carrying
#      the prior "2022-Q3" totals, pushing "2021-Q3" down one row.

$wb = $excel.ActiveWorkbook

$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# 1. Duplicate the "2022-Q3" sheet; the copy is inserted immediately
#    after it and keeps the old "2022-Q3" data untouched.
$wsQ3.Copy($null, $wsQ3)
$wsOldQ3Copy = $wb.Worksheets.Item($wsQ3.Index + 1)

# 2. The original sheet becomes "2022-Q4" with refreshed figures.
$wsQ3.Name = "2022-Q4"
$wsQ4 = $wsQ3

# These figures are stored as text (matching the source sheet's
# original inline-string cells), so force text entry and then drop
# back to the "Normal" style to avoid leaving a stray number-format
# override behind.
foreach ($pair in @(@("D2","4.06"), @("E2","93.21"), @("F2","3.11"), @("G2","0.1263"))) {
    $cell = $wsQ4.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

# 3. The duplicate keeps the old numbers and is renamed back to "2022-Q3".
$wsOldQ3Copy.Name = "2022-Q3"

# 4. Update the "总计" summary sheet: row 2 -> 2022-Q4 (new totals), a new
#    row 3 -> 2022-Q3 (the totals row 2 used to hold), row 4 -> 2021-Q3
#    (previously row 3, values unchanged).
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.13

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.11

# Row 4 is brand new, so it doesn't inherit the "index column" style
# that A2/A3 already carry. Copy that formatting across before writing
# the value.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.08
